$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Insert a new row before row 92, pushing existing rows 92-102 down to 93-103
$ws.Rows.Item(92).Insert()

# Populate the new row 92 with the e061 entry
# (the row Insert above already copied column A/B's formatting down, so the
# new cells keep the same styles as the rest of the table)
$ws.Range("A92").Value = "e061"
$ws.Range("B92").Value = "<Bold>e061 Crew Switch</Bold> " + [char]10 + "<InlineUIContainer><Button Content='r19.22' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   " + [char]10 + "<LineBreak/><LineBreak/>" + [char]10 + "The assistant driver moves through the tank ro replace the incapacitated crewman. The assistant driver takes on the role but at half rating. Click image to  continue." + [char]10 + "<LineBreak/><LineBreak/>" + [char]10 + "                                            <InlineUIContainer><Image Name='CarryingMan' Height='80' Width='200'></Image></InlineUIContainer>"

$ws.Rows.Item(92).RowHeight = 90

# Update the visible view / selection to match the authored state
$excel.ActiveWindow.ScrollRow = 90
$ws.Range("B93").Select()
